# Update the date heading
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-01-31 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-02-01 Saturday", 2) | Out-Null

# Simple 1:1 value replacements inside the multiplication table.
# Each "old" string is unique in the document, so Find/Execute is safe here.
$replacements = @(
    @("95×60=5700", "57×89=5073"),
    @("40×81=3240", "50×51=2550"),
    @("99×28=2772", "13×42=546"),
    @("27×77=2079", "18×84=1512"),
    @("78×83=6474", "42×87=3654"),
    @("53×29=1537", "51×74=3774"),
    @("51×28=1428", "56×95=5320"),
    @("58×27=1566", "16×87=1392"),
    @("50×77=3850", "47×43=2021"),
    @("59×70=4130", "14×51=714"),
    @("27×51=1377", "37×66=2442"),
    @("71×12=852",  "30×23=690"),
    @("69×20=1380", "73×63=4599"),
    @("25×43=1075", "41×77=3157"),
    @("88×32=2816", "79×32=2528"),
    @("91×70=6370", "98×55=5390"),
    @("76×71=5396", "31×17=527"),
    @("19×20=380",  "97×31=3007"),
    @("13×85=1105", "71×53=3763"),
    @("25×86=2150", "23×76=1748")
)

foreach ($pair in $replacements) {
    $d.Content.Find.Execute($pair[0], $true, $false, $false, $false, $false,
                             $true, 1, $false, $pair[1], 2) | Out-Null
}

# The last table row is reshuffled: the first cell's value moves on, three
# brand-new problems are inserted after it, and the old trailing three cells
# are dropped, while the table keeps its 5 columns. Address the row's cells
# directly by position so the duplicate-valued "68×48=3264" text can't be
# confused with the other Find/Replace operations above.
$t = $d.Tables.Item(1)
$lastRow = $t.Rows.Count

$lastRowValues = @("31×28=868", "40×25=1000", "92×31=2852", "67×20=1340", "68×48=3264")
for ($col = 1; $col -le $lastRowValues.Count; $col++) {
    $t.Cell($lastRow, $col).Range.Text = $lastRowValues[$col - 1]
}
